# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextCell 'D2' '306.33'
Set-TextCell 'D3' '41.24'
Set-TextCell 'E3' '5.20%'
Set-TextCell 'D4' '5.105'
Set-TextCell 'E4' '2.46%'
Set-TextCell 'D5' '0.07606'
Set-TextCell 'E5' '-1.45%'
Set-TextCell 'B6' 'FTXToken'
Set-TextCell 'C6' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell 'D6' '1.619'
Set-TextCell 'E6' '1.95%'
Set-TextCell 'B7' 'BTSEToken'
Set-TextCell 'C7' 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell 'D7' '2.454'
Set-TextCell 'E7' '-5.28%'
Set-TextCell 'B8' 'MXToken'
Set-TextCell 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 'D8' '0.9062'
Set-TextCell 'E8' '-0.41%'
Set-TextCell 'B9' 'LiechtensteinCryptoassetsExchange'
Set-TextCell 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextCell 'D9' '0.1019'
Set-TextCell 'E9' '0.70%'
Set-TextCell 'B10' 'WazirX'
Set-TextCell 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextCell 'D10' '0.1752'
Set-TextCell 'E10' '1.54%'
Set-TextCell 'B11' 'MandalaExchangeToken'
Set-TextCell 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextCell 'D11' '0.09091'
Set-TextCell 'E11' '1.05%'
Set-TextCell 'B12' 'BitrueCoin'
Set-TextCell 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextCell 'D12' '0.04263'
Set-TextCell 'E12' '-4.05%'
Set-TextCell 'B13' 'BitMartToken'
Set-TextCell 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell 'D13' '0.1056'
Set-TextCell 'E13' '-0.18%'
Set-TextCell 'B14' 'BitForexToken'
Set-TextCell 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell 'D14' '0.001228'
Set-TextCell 'E14' '-3.59%'
Set-TextCell 'B15' 'TigerCash'
Set-TextCell 'C15' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell 'D15' '0.005859'
Set-TextCell 'E15' '3.83%'
Set-TextCell 'B16' 'LEO'
Set-TextCell 'C16' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 'D16' '3.350'
Set-TextCell 'E16' '-0.30%'
Set-TextCell 'B17' 'GateToken'
Set-TextCell 'C17' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextCell 'D17' '4.267'
Set-TextCell 'E17' '-0.04%'
Set-TextCell 'D18' '0.3274'
Set-TextCell 'E18' '-2.75%'
Set-TextCell 'D19' '6.551'
Set-TextCell 'E19' '-7.04%'
Set-TextCell 'D20' '0.1357'
Set-TextCell 'E20' '0.92%'
Set-TextCell 'D21' '0.2725'
Set-TextCell 'D22' '0.04181'
Set-TextCell 'E22' '1.03%'
Set-TextCell 'D23' '0.001228'
Set-TextCell 'E23' '3.02%'
Set-TextCell 'D24' '0.004077'
Set-TextCell 'E24' '-0.02%'
Set-TextCell 'E25' '6.32%'
Set-TextCell 'D26' '0.0003009'
Set-TextCell 'E26' '0.75%'
Set-TextCell 'D38' '0.02376'
Set-TextCell 'E38' '1.37%'
Set-TextCell 'D39' '0.05144'
Set-TextCell 'E39' '0.49%'
Set-TextCell 'D40' '0.007780'
Set-TextCell 'E40' '-2.50%'
Set-TextCell 'D41' '0.1296'
Set-TextCell 'E41' '-2.19%'
Set-TextCell 'D42' '0.006990'
Set-TextCell 'E42' '-7.92%'
Set-TextCell 'E43' '-4.33%'
Set-TextCell 'D44' '0.008454'
Set-TextCell 'E44' '5.80%'
Set-TextCell 'D45' '0.3332'
Set-TextCell 'E45' '0.50%'
Set-TextCell 'D46' '0.00006357'
Set-TextCell 'E46' '-5.03%'
Set-TextCell 'E47' '-0.24%'
Set-TextCell 'D48' '0.009025'
Set-TextCell 'E48' '164.46%'
Set-TextCell 'D49' '0.004404'
Set-TextCell 'E49' '7.02%'
Set-TextCell 'D50' '0.00002101'
Set-TextCell 'E50' '-0.24%'
Set-TextCell 'D51' '0.0002001'
Set-TextCell 'E51' '-0.24%'
